$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new log entry to the first empty row following the existing data in column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "2025-11-24 13:12:12"
